# Added some test case
# Reorders the data rows (rows 3-21, columns A-F) of the worksheet.
# The set of row-tuples is unchanged; only their row assignment is shuffled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(701, 3, 90, 45, 97, 15),
    @(901, 16, 15, 45, 60, 60),
    @(902, 1, 0, 0, 0, 0),
    @(501, 9, 52, 30, 75, 45),
    @(801, 3, 67, 65, 52, 45),
    @(301, 6, 45, 30, 60, 45),
    @(601, 9, 60, 67, 60, 42),
    @(1202, 2, 10, 10, 10, 10),
    @(101, 9, 30, 15, 60, 15),
    @(401, 9, 48, 67, 75, 45),
    @(201, 9, 30, 15, 45, 30),
    @(1201, 2, 10, 10, 10, 10),
    @(1203, 3, 15, 15, 15, 15),
    @(802, 0, 4, 5, 4, 0),
    @(1101, 0, 15, 30, 30, 0),
    @(3, 0, 3, 3, 3, 3),
    @(502, 0, 4, 0, 0, 0),
    @(1, 0, 2, 2, 2, 2),
    @(2, 0, 2, 2, 2, 2)
)

$startRow = 3
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 0; $col -lt $values.Length; $col++) {
        $ws.Cells.Item($row, $col + 1).Value = $values[$col]
    }
}
